# Update file upload functionality
# Appends one new logged data row to each of the four worksheets, mirroring
# the structure/format of the existing rows already present on each sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: ROW50-FE-LIFTER -> append row 46 -----------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A46").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Range("A46").Value2 = 45747.18384898148
$ws1.Range("B46").Value = "0x01,0x90"
$ws1.Range("C46").Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
$ws1.Range("D46").Value = "0x01,0x6a"
$ws1.Range("E46").Value = "0xe"
$ws1.Range("F46").Value2 = 400
$ws1.Range("G46").Value2 = 568631262647114000000000.0
$ws1.Range("H46").Value2 = 362
$ws1.Range("I46").Value2 = 14

# --- Sheet 2: ROW50-MID-LIFTER -> append row 48 -----------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A48").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Range("A48").Value2 = 45747.1524074074
$ws2.Range("B48").Value = "0x01,0x90 "
$ws2.Range("C48").Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws2.Range("D48").Value = "0x01,0x6e"
$ws2.Range("E48").Value = "0x19"
$ws2.Range("F48").Value2 = 400
# This particular column stores a 24-digit figure that exceeds normal
# numeric precision as literal text (matches the rest of column G on this
# sheet), so force it in as text rather than letting it be parsed as a number.
$ws2.Range("G48").Value = "'568631262647113771663628"
$ws2.Range("H48").Value2 = 366
$ws2.Range("I48").Value2 = 25

# --- Sheet 3: ROW11-FE-LIFTER -> append row 46 ------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("A46").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws3.Range("A46").Value2 = 45747.206095625
$ws3.Range("B46").Value = "0x01,0x90"
$ws3.Range("C46").Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
$ws3.Range("D46").Value = "0x01,0x6a"
$ws3.Range("E46").Value = "0x14"
$ws3.Range("F46").Value2 = 400
$ws3.Range("G46").Value2 = 568631262647114000000000.0
$ws3.Range("H46").Value2 = 362
$ws3.Range("I46").Value2 = 20

# --- Sheet 4: ROW11-MID-LIFTER -> append row 46 -----------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("A46").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws4.Range("A46").Value2 = 45747.34497392361
$ws4.Range("B46").Value = "0x01,0x90"
$ws4.Range("C46").Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
$ws4.Range("D46").Value = "0x01,0x6e"
$ws4.Range("E46").Value = "0x19"
$ws4.Range("F46").Value2 = 400
$ws4.Range("G46").Value2 = 568631262647114000000000.0
$ws4.Range("H46").Value2 = 366
$ws4.Range("I46").Value2 = 25
